$d = $word.ActiveDocument

# --- Step 1: find the paragraph that follows "Linear interpolation necessary
#     to evaluate models at specific times" (an empty paragraph) -- the new
#     webinar-notes paragraphs get appended right after it.
$anchor = $d.Content
$anchor.Find.Execute("Linear interpolation necessary to evaluate models at specific times", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorIndex = $anchor.Paragraphs.First.Index
$emptyParagraph = $d.Paragraphs.Item($anchorIndex + 1)
$insertAt = $emptyParagraph.Range.Duplicate
$insertAt.Collapse(0)
$insertAt.InsertParagraphAfter()
$newParagraphRange = $d.Paragraphs.Item($anchorIndex + 2).Range

# --- Step 2: mint a fresh hyperlink relationship for the youtube link (the
#     same relationship Word would register) by adding it to a scratch
#     paragraph at the very end of the document, then deleting that scratch
#     paragraph again. The relationship entry survives in
#     word/_rels/document.xml.rels even though the visible run is gone, so
#     we can reference its id explicitly in the OOXML we insert below.
$scratchAnchor = $d.Content
$scratchAnchor.Collapse(0)
$scratchAnchor.InsertParagraphAfter()
$scratchParagraph = $d.Paragraphs.Last
$hyperlink = $d.Hyperlinks.Add($scratchParagraph.Range, "https://www.youtube.com/watch?v=kWlaGmsh9Mg")
$cleanupRange = $d.Range($scratchParagraph.Range.Start, $d.Content.End)
$cleanupRange.Delete()

# --- Step 3: insert the new paragraphs (hyperlink + five bullet notes) as
#     literal OOXML right after the empty paragraph located in step 1,
#     referencing the relationship minted in step 2 (rId13).
$xmlPayload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:rPr/></w:pPr><w:hyperlink r:id="rId13"><w:r><w:rPr><w:color w:val="1155cc"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">https://www.youtube.com/watch?v=kWlaGmsh9Mg</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Newton’s 2nd law allows future positions and velocities of objects to be predicted if current positions and velocities are known</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Quantum mechanics introduces randomness at smallest scales - impossible to measure position and momentum with infinite precision</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Simple deterministic systems with few elements exhibit seemingly random behaviour - chaos - due to non-linearity and sensitivity to initial conditions</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Interconnectedness of universe means causality is not easily computable</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Chaos = unpredictable oscillations</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newParagraphRange.InsertXML($xmlPayload)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
